$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Text used for the new "years back" labels (replacing "2yr" / "all (4 years)"
# / the bare "2" code) across the three repeated blocks (B:E, G:J, L:O).
# ---------------------------------------------------------------------------
$oneYear  = "one prior year of data (2001–2002)"
$twoYear  = "two prior years of data (2000 to 2002)"
$allYears = "all five years that are available in the CMS panel (1997 to 2002)"

$modelCols = @("B", "G", "L")
$yearsCols = @("C", "H", "M")

for ($i = 0; $i -lt $modelCols.Length; $i++) {
    $mCol = $modelCols[$i]
    $yCol = $yearsCols[$i]

    # Row 4-6 => model 1 ("Calculate VAM..." wiped out, replaced by plain 1)
    # Row 7-9 => model 2 ("Add also only..." wiped out, replaced by plain 2)
    $ws.Range($mCol + "4").Value = 1
    $ws.Range($mCol + "5").Value = 1
    $ws.Range($mCol + "6").Value = 1
    $ws.Range($mCol + "7").Value = 2
    $ws.Range($mCol + "8").Value = 2
    $ws.Range($mCol + "9").Value = 2

    # "years back" column now holds descriptive text in every row instead of
    # a mix of numbers/"2yr"/"all (4 years)".
    $ws.Range($yCol + "4").Value = $oneYear
    $ws.Range($yCol + "5").Value = $twoYear
    $ws.Range($yCol + "6").Value = $allYears
    $ws.Range($yCol + "7").Value = $oneYear
    $ws.Range($yCol + "8").Value = $twoYear
    $ws.Range($yCol + "9").Value = $allYears

    # The old rows 7-9 carried a wrap-text style (tall 156.75pt rows) for the
    # long "Add also only a third-order polynomial..." text that no longer
    # exists. Put those cells back to the plain Normal style and let Excel
    # re-measure the row height now that nothing needs wrapping.
    $ws.Range($mCol + "7:" + $yCol + "9").Style = "Normal"
}

$ws.Range("7:9").EntireRow.AutoFit()

# New trailing marker cell.
$ws.Range("C12").Value = 1

# Column L's custom "best fit" width (sized for the old long text) is no
# longer needed; column M now needs a width instead.
$ws.Columns.Item(13).ColumnWidth = 17.17

$ws.Range("C4:C9").Select()
